# Apply updated cryptocurrency price/volume figures to Sheet1.
# Generated from the authoritative cell-level diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.348.16'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '1.709.67'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5292'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06640'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.86'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07683'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.509'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.18%  '
$ws.Range('D13').Value = '1.945.02'
$ws.Range('E13').Value = '  -0.98%  '
$ws.Range('D14').Value = '1.706.40'
$ws.Range('E14').Value = '  -1.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5821'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('D16').Value = '0.0₅8187'
$ws.Range('E16').Value = '  -1.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.94'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').Value = '27.328.18'
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '220.98'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.634'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.93%  '
$ws.Range('E22').Value = '  -1.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.986'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.73%  '
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.90%  '
$ws.Range('E26').Value = '  -1.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1210'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.257'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.28'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05365'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.295'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.460'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.428'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.646'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.861'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9520'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.399'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5865'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01638'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Value = '1.082.42'
$ws.Range('E40').Value = '  +2.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.805'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.13%  '
$ws.Range('E42').Value = '  -1.26%  '
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('D45').Value = '1.852.42'
$ws.Range('E46').Value = '  +2.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '57.85'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4522'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.86%  '
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.093'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05230'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.52%  '
